$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.189.77"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "2.376.32"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.48%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("D9").Value = "2.377.90"
$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.01%  "

$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("E12").Value = "  +0.69%  "

$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000167"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.787.44"
$ws.Range("E16").Value = "  -1.23%  "

$ws.Range("D17").Value = "61.109.34"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "2.389.83"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.19%  "

$ws.Range("E25").Value = "  -13.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.11%  "

$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").Value = "2.487.19"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.54%  "

$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "506.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.61%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0888"
$ws.Range("E31").Value = "  -4.98%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.150"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.81%  "

$ws.Range("E33").Value = "  -4.13%  "

$ws.Range("E34").Value = "  -0.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.66"
$ws.Range("D37").Style = "Normal"

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.379"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "146.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.99%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.20%  "

$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.06%  "

$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("E49").Value = "  -4.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.576"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("E51").Value = "  +0.34%  "
